# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K" - strikeouts) values for rows 2-45 on Sheet1 to the
# newly recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 3
    9  = 2
    10 = 0
    11 = 3
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 3
    17 = 2
    18 = 3
    19 = 1
    20 = 2
    21 = 1
    22 = 2
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 0
    28 = 3
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 2
    41 = 1
    42 = 0
    43 = 2
    44 = 3
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
